{"js": "const body = context.document.body;\n\nconst fr6Old = \"FR6: The system shall interpret the first two digits of a BasicML word as the opcode.\";\nconst fr6New = \"FR6: The system shall interpret the first three digits of a BasicML word as the opcode.\";\nconst fr7Old = \"FR7: The system shall interpret the last two digits of a BasicML word as the operand.\";\nconst fr7New = \"FR7: The system shall interpret the last three digits of a BasicML word as the operand.\";\n\nconst fr6Results = body.search(fr6Old, { matchCase: true });\nfr6Results.load(\"items\");\nawait context.sync();\nif (fr6Results.items.length > 0) {\n  fr6Results.items[0].insertText(fr6New, \"Replace\");\n}\n\nconst fr7Results = body.search(fr7Old, { matchCase: true });\nfr7Results.load(\"items\");\nawait context.sync();\nif (fr7Results.items.length > 0) {\n  fr7Results.items[0].insertText(fr7New, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"FR6: The system shall interpret the first two digits of a BasicML word as the opcode.\"\n        New = \"FR6: The system shall interpret the first three digits of a BasicML word as the opcode.\"\n    },\n    @{\n        Old = \"FR7: The system shall interpret the last two digits of a BasicML word as the operand.\"\n        New = \"FR7: The system shall interpret the last three digits of a BasicML word as the operand.\"\n    }\n)\n\n$count = $d.Paragraphs.Count\nforeach ($rep in $replacements) {\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        $trimmed = $t.TrimEnd([char]13, [char]7)\n        if ($trimmed -eq $rep.Old) {\n            $p.Range.Text = $rep.New\n            break\n        }\n    }\n}\n\nWrite-Output \"done\"\n"}
